$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.465275764465332
$ws.Range("B1").Value = 2.658677577972412
$ws.Range("C1").Value = 2.006633520126343
$ws.Range("D1").Value = 1.900489687919617
$ws.Range("E1").Value = 2.009080410003662
